$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LMSProd")
$ws.Range("A2").Value = "test"
